$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 279. Excel shifts rows
# 279..290 down to 281..292, preserving all their existing values/styles.
$ws.Rows.Item(279).Insert()
$ws.Rows.Item(279).Insert()

# New row 279 (week of 44509, "Primera")
$ws.Range("A279").Value = 9
$ws.Range("B279").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C279").Value = "Metropolitana"
$ws.Range("D279").Value = 44509
$ws.Range("E279").Value = 13
$ws.Range("F279").Value = 100114014
$ws.Range("G279").Value = "Betarraga"
$ws.Range("H279").Value = "Sin especificar"
$ws.Range("I279").Value = "Primera"
$ws.Range("J279").Value = 7900
$ws.Range("K279").Value = 90
$ws.Range("L279").Value = 100
$ws.Range("M279").Value = 95
$ws.Range("N279").Value = '$/unidad'
$ws.Range("O279").Value = "Región Metropolitana"
$ws.Range("P279").Value = 95
$ws.Range("Q279").Value = 1
$ws.Range("R279").Value = "Hortaliza"

# New row 280 (week of 44509, "Segunda")
$ws.Range("A280").Value = 9
$ws.Range("B280").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C280").Value = "Metropolitana"
$ws.Range("D280").Value = 44509
$ws.Range("E280").Value = 13
$ws.Range("F280").Value = 100114014
$ws.Range("G280").Value = "Betarraga"
$ws.Range("H280").Value = "Sin especificar"
$ws.Range("I280").Value = "Segunda"
$ws.Range("J280").Value = 3400
$ws.Range("K280").Value = 60
$ws.Range("L280").Value = 70
$ws.Range("M280").Value = 65
$ws.Range("N280").Value = '$/unidad'
$ws.Range("O280").Value = "Región Metropolitana"
$ws.Range("P280").Value = 65
$ws.Range("Q280").Value = 1
$ws.Range("R280").Value = "Hortaliza"
